$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 01:52"

# Row 4 - Estados Unidos: updated case counts
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 336327
$ws.Range("C4").Value = 24970
$ws.Range("D4").Value = 17245
$ws.Range("E4").Value = 309477
$ws.Range("F4").Value = 8702
$ws.Range("G4").Value = 1154
$ws.Range("H4").Value = 9605

# Rows 23-24 - Australia moves above Noruega with fresh data; Noruega keeps its previous data
$ws.Range("A23").Value = "Australia"
$ws.Range("B23").Value = 5750
$ws.Range("C23").Value = 200
$ws.Range("D23").Value = 2315
$ws.Range("E23").Value = 3398
$ws.Range("F23").Value = 91
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 37

$ws.Range("A24").Value = "Noruega"
$ws.Range("B24").Value = 5687
$ws.Range("C24").Value = 137
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = 5584
$ws.Range("F24").Value = 89
$ws.Range("G24").Value = 9
$ws.Range("H24").Value = 71

# Rows 30-34 - India moves above Polonia with fresh data; Polonia, Rumania, Malasia and
# Ecuador each shift down one row keeping their own previous data
$ws.Range("A30").Value = "India"
$ws.Range("B30").Value = 4288
$ws.Range("C30").Value = 700
$ws.Range("D30").Value = 328
$ws.Range("E30").Value = 3843
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 18
$ws.Range("H30").Value = 117

$ws.Range("A31").Value = "Polonia"
$ws.Range("B31").Value = 4102
$ws.Range("C31").Value = 475
$ws.Range("D31").Value = 134
$ws.Range("E31").Value = 3874
$ws.Range("F31").Value = 50
$ws.Range("G31").Value = 15
$ws.Range("H31").Value = 94

$ws.Range("A32").Value = "Rumania"
$ws.Range("B32").Value = 3864
$ws.Range("C32").Value = 251
$ws.Range("D32").Value = 374
$ws.Range("E32").Value = 3339
$ws.Range("F32").Value = 141
$ws.Range("G32").Value = 5
$ws.Range("H32").Value = 151

$ws.Range("A33").Value = "Malasia"
$ws.Range("B33").Value = 3662
$ws.Range("C33").Value = 179
$ws.Range("D33").Value = 1005
$ws.Range("E33").Value = 2596
$ws.Range("F33").Value = 99
$ws.Range("G33").Value = 4
$ws.Range("H33").Value = 61

$ws.Range("A34").Value = "Ecuador"
$ws.Range("B34").Value = 3646
$ws.Range("C34").Value = 181
$ws.Range("D34").Value = 100
$ws.Range("E34").Value = 3366
$ws.Range("F34").Value = 100
$ws.Range("G34").Value = 8
$ws.Range("H34").Value = 180

# Rows 132-135 - Guatemala and Guayana Francesa move above Aruba and El Salvador with
# fresh data; Aruba and El Salvador each shift down one row keeping their own previous data
$ws.Range("A132").Value = "Guatemala"
$ws.Range("B132").Value = 70
$ws.Range("C132").Value = 9
$ws.Range("D132").Value = 15
$ws.Range("E132").Value = 52
$ws.Range("F132").Value = 3
$ws.Range("G132").Value = 1
$ws.Range("H132").Value = 3

$ws.Range("A133").Value = "Guayana Francesa"
$ws.Range("B133").Value = 68
$ws.Range("C133").Value = 7
$ws.Range("D133").Value = 27
$ws.Range("E133").Value = 41
$ws.Range("F133").Value = 1
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 0

$ws.Range("A134").Value = "Aruba"
$ws.Range("B134").Value = 64
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 1
$ws.Range("E134").Value = 63
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 0

$ws.Range("A135").Value = "El Salvador"
$ws.Range("B135").Value = 62
$ws.Range("C135").Value = 6
$ws.Range("D135").Value = 2
$ws.Range("E135").Value = 57
$ws.Range("F135").Value = 4
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 3
